$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S (year 2022) values, one per row, mirroring existing column R formatting.
$values = @{
    4  = 2022
    5  = 4.9538761752705343
    6  = 11.304954640614097
    7  = 5.1593323216995444
    8  = 13.687943262411348
    9  = 10.22864019253911
    10 = 9.1213700670141478
    11 = 3.1335149863760217
    12 = 2.872905173311127
    13 = 3.527842284697861
    14 = 5.0305321314335565
}

# xlPasteFormats = -4122 : copy the formatting of column R into the new column S
# so each new cell matches the style of the corresponding cell to its left.
for ($row = 4; $row -le 14; $row++) {
    $srcCell = $ws.Range("R$row")
    $dstCell = $ws.Range("S$row")
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null
    $dstCell.Value = $values[$row]
}

$excel.CutCopyMode = 0

# Update the active selection, as recorded in the workbook view.
$ws.Range("T6").Select() | Out-Null
